$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.463.04"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "1.838.66"
$ws.Range("E3").Value = "  +1.31%  "

$ws.Range("E4").Value = "  +1.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.02"
$ws.Range("E5").Value = "  +1.81%  "

$ws.Range("E6").Value = "  +0.95%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").Value = "  +1.64%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3691"
$ws.Range("E8").Value = "  +0.85%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").Value = "  +1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8856"
$ws.Range("E10").Value = "  +1.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.45"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("D12").Value = "1.926.21"
$ws.Range("E12").Value = "  +5.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07336"
$ws.Range("E13").Value = "  +3.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.452"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.28"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.582"
$ws.Range("E16").Value = "  +1.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.011"
$ws.Range("E17").Value = "  +0.78%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008820"
$ws.Range("E18").Value = "  +1.36%  "

$ws.Range("E19").Value = "  +0.99%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.81"
$ws.Range("E20").Value = "  +1.05%  "

$ws.Range("D21").Value = "27.495.87"
$ws.Range("E21").Value = "  +2.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.326"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("E23").Value = "  +0.56%  "

$ws.Range("D24").Value = "2.134.59"
$ws.Range("E24").Value = "  +3.60%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.907"
$ws.Range("E25").Value = "  +0.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.09"
$ws.Range("E26").Value = "  +0.68%  "

$ws.Range("E27").Value = "  +1.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.146"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.251"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.96"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09001"
$ws.Range("E31").Value = "  +1.16%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7560"
$ws.Range("E32").Value = "  +0.27%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.183"
$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.560"
$ws.Range("E34").Value = "  +1.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.953"
$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.013"
$ws.Range("E36").Value = "  +1.12%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.104"
$ws.Range("E37").Value = "  +1.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05334"
$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01956"
$ws.Range("E39").Value = "  +0.41%  "

$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.994"
$ws.Range("E40").Value = "  +0.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.336"
$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.406"
$ws.Range("E42").Value = "  +5.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5337"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1661"
$ws.Range("E44").Value = "  +0.49%  "

$ws.Range("E45").Value = "  +0.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4912"
$ws.Range("E46").Value = "  +0.78%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.53"
$ws.Range("E47").Value = "  +1.52%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "105.22"
$ws.Range("E48").Value = "  +1.95%  "

$ws.Range("E49").Value = "  +1.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.679"
$ws.Range("E50").Value = "  +1.14%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06322"
$ws.Range("E51").Value = "  +0.49%  "
